# Apply "More test results added" edit:
#  - fill in previously-empty G/H result cells for rows 15, 16, 18, 21, 24
#  - move the active selection on Sheet1 from J30 to G14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G15").Value = 0.37605019996408301
$ws.Range("H15").Value = 5.57825930003309

$ws.Range("G16").Value = 0.43474626481367001
$ws.Range("H16").Value = 5.74009769930543

$ws.Range("G18").Value = 0.41064017314810902
$ws.Range("H18").Value = 0.36200941342556398

$ws.Range("G21").Value = 0.42678876944356198
$ws.Range("H21").Value = 0.24662390833381501

$ws.Range("G24").Value = 0.38344519444055702
$ws.Range("H24").Value = 0.36089901789855

# Update the selected cell/range shown when the sheet is reopened.
$null = $ws.Range("G14").Select()
